# ------------------------------------------------------------------
# Add new "Cleaning" sheet translations (Brewer / Milker cleaning
# strings) in English (col B) and Italian (col C), with their
# $LAB_ label codes (col A). Written in the exact order the rows
# were authored so the shared-string table lines up with the source
# workbook.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cleaning")

$ws.Range("B14").Value = 'Brewer Cleaning is not started (or ended)'
$ws.Range("A14").Value = '$LAB_CLEANSAN_1'
$ws.Range("B15").Value = 'Brewer Cleaning is started'
$ws.Range("A15").Value = '$LAB_CLEANSAN_2'
$ws.Range("B16").Value = 'Brewer placed'
$ws.Range("B17").Value = 'Put pastille and push START'
$ws.Range("B18").Value = 'Infusion'
$ws.Range("B19").Value = 'Brewer cleaning cycles'
$ws.Range("A16").Value = '$LAB_CLEANSAN_3'
$ws.Range("A17").Value = '$LAB_CLEANSAN_4'
$ws.Range("A18").Value = '$LAB_CLEANSAN_5'
$ws.Range("A19").Value = '$LAB_CLEANSAN_6'
$ws.Range("B20").Value = 'Repeat cleaning ?'
$ws.Range("A20").Value = '$LAB_CLEANSAN_7'
$ws.Range("B21").Value = 'Brewer placed in brush position, press CONTINUE when finished'
$ws.Range("A21").Value = '$LAB_CLEANSAN_8'
$ws.Range("B22").Value = 'Skip final coffee or make a coffee'
$ws.Range("B23").Value = 'Coffee delivery'
$ws.Range("B24").Value = 'Rinsing'
$ws.Range("A22").Value = '$LAB_CLEANSAN_9'
$ws.Range("A23").Value = '$LAB_CLEANSAN_10'
$ws.Range("A24").Value = '$LAB_CLEANSAN_11'
$ws.Range("B26").Value = 'Milker Cleaning is started'
$ws.Range("A26").Value = '$LAB_CLEANMILK_1'
$ws.Range("B28").Value = 'Wait for confirm'
$ws.Range("B31").Value = 'Wait for second confirm'
$ws.Range("B30").Value = 'Warming for water'
$ws.Range("B27").Value = 'Warming for cleaner'
$ws.Range("A27").Value = '$LAB_CLEANMILK_2'
$ws.Range("A28").Value = '$LAB_CLEANMILK_3'
$ws.Range("A29").Value = '$LAB_CLEANMILK_4'
$ws.Range("A30").Value = '$LAB_CLEANMILK_5'
$ws.Range("A31").Value = '$LAB_CLEANMILK_6'
$ws.Range("B29").Value = 'Doing cleaner cycles (12)'
$ws.Range("C14").Value = 'Lavaggio non iniziato (o terminato)'
$ws.Range("C15").Value = 'Lavaggio del gruppo in corso'
$ws.Range("C16").Value = 'Gruppo posizionato'
$ws.Range("C17").Value = 'Inserire la pastiglia e premere INIZIA'
$ws.Range("C18").Value = 'Infusione'
$ws.Range("C19").Value = 'Ciclo di pulizia del gruppo'
$ws.Range("C20").Value = 'Ripetere la pulizia?'
$ws.Range("C21").Value = 'Gruppo posizionato, premere CONTINUA quando terminato'
$ws.Range("C22").Value = 'Si desidera erogare un caffè?'
$ws.Range("C23").Value = 'Erogazione in corso'
$ws.Range("C24").Value = 'Risciacquo'
$ws.Range("C26").Value = 'Lavaggio del cappuccinatore in corso'
$ws.Range("C27").Value = 'Attesa temperatura'
$ws.Range("C28").Value = 'Attesa conferma'
$ws.Range("C29").Value = 'Cicli di lavaggio in corso (12)'
$ws.Range("C30").Value = 'Attesa temperatura'
$ws.Range("C31").Value = 'Attesa seconda conferma'

# Widen columns B (English) and C (Italian) to fit the new, longer
# strings.
$ws.Columns.Item(2).ColumnWidth = 58.5
$ws.Columns.Item(3).ColumnWidth = 23.833333333333332

# Move the active window / selection from "Devices" (previously the
# active tab) onto "Cleaning", scrolled to the newly added rows with
# C32 selected -- matches the saved view state in the workbook.
$ws.Activate()
$ws.Range("C32").Select()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
